$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2021")

$ws.Range("B2").Value = 0.1183284763020911
$ws.Range("C2").Value = 0.6112895150387458
$ws.Range("D2").Value = 0.5421014366391211
$ws.Range("E2").Value = 0.7362753809812747
$ws.Range("F2").Value = 0.7477730290289269
$ws.Range("G2").Value = 18
